# Update "想去人数" (interested count) values in column F for three events,
# on both the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# Sheet "展览": rows 4, 10, 13
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 3641
$wsExhibit.Range("F10").Value = 73
$wsExhibit.Range("F13").Value = 1995

# Sheet "全部类型": rows 4, 11, 16
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 3641
$wsAll.Range("F11").Value = 73
$wsAll.Range("F16").Value = 1995
